$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "1.005") are not
# auto-converted to numbers by Excel, matching the original inline-string semantics.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.675.65"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.696.01"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").Value = "315.64"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "0.3921"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "1.506"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "53.09"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "0.08769"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "7.675"
$ws.Range("E13").Value = "  +5.29%  "
$ws.Range("D14").Value = "24.44"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").Value = "0.00001362"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "8.000"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "1.694.72"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "98.53"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "0.07123"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "7.392"
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("D22").Value = "1.009"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "24.641.39"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "3.032"
$ws.Range("E25").Value = "  -7.35%  "
$ws.Range("D26").Value = "2.356"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "22.75"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "162.54"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "8.489"
$ws.Range("E29").Value = "  +13.19%  "
$ws.Range("D30").Value = "137.85"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").Value = "5.237"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "1.884.73"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "0.08884"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").Value = "7.498"
$ws.Range("E34").Value = "  +4.79%  "
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").Value = "1.995"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "0.02930"
$ws.Range("E37").Value = "  +7.18%  "
$ws.Range("D38").Value = "0.2736"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -5.52%  "
$ws.Range("D40").Value = "14.34"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Value = "0.09146"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "0.7900"
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "16.60"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").Value = "0.7230"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").Value = "2.575"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").Value = "4.225"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D49").Value = "1.333"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").Value = "139.34"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "91.65"
$ws.Range("E51").Value = "  +1.73%  "
